$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.259.36"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "2.024.45"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'248.05"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "'58.00"
$ws.Range("E7").Value = "  -3.50%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.388"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").Value = "'0.0801"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'14.90"
$ws.Range("E12").Value = "  +4.70%  "
$ws.Range("D13").Value = "2.322.74"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "'0.830"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "'5.37"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "2.025.24"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("D18").Value = "37.190.44"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "'69.94"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").Value = "'228.53"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'2.55"
$ws.Range("E24").Value = "  +4.46%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'163.46"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  -5.84%  "
$ws.Range("D29").Value = "'19.89"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'4.76"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "'0.0666"
$ws.Range("E33").Value = "  +8.15%  "
$ws.Range("D34").Value = "'4.55"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  +8.48%  "
$ws.Range("E36").Value = "  +6.15%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").Value = "'5.34"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("D41").Value = "'0.0969"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "'1.16"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'16.41"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").Value = "1.400.51"
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("D46").Value = "'90.93"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("E47").Value = "  +4.03%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").Value = "'2.07"
$ws.Range("E49").Value = "  +11.25%  "
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "2.214.70"
$ws.Range("E51").Value = "  +2.69%  "
